$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "20.006.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.54%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.420.49"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.04%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -1.22%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9997"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.82%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "277.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.72%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3692"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.10%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3102"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.89%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.83"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.61%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.047"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.33%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06553"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.02%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9997"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.35%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.519"
$ws.Range("D13").Style = "Normal"

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.96%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.223"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.67%  "

# Row 16
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.420.70"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.83%  "

# Row 17
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001024"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.70%  "

# Row 18
$ws.Range("E18").Value = "  -12.50%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9997"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.82%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -10.94%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.629"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.00%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.79"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.30%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.09%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.238"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.20%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "20.013.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.53%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.293"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.72%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "133.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -8.36%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.37"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.67%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.581.70"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.82%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "110.39"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.58%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.906"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -18.05%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.261"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.83%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8233"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -9.27%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07763"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.31%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.488"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.08%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.291"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.80%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.924"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.64%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05857"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.57%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9993"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.82%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02067"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.64%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.53"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.54%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1890"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.42%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.104"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.22%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5337"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.23%  "

# Row 45
$ws.Range("E45").Value = "  -1.39%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.547"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.07%  "

# Row 47
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "117.31"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.68%  "

# Row 48
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5213"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.93%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.778"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.50%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.037"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.03%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.000"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.83%  "
